# Auto-generated edit script applying scraped-data updates to Kraken_Profits workbook
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H,I,J,K,L,M,N)
# across the ALC, ARM, CRP, GSM, LTW and WVR sheets, matching a refreshed price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 25199.5
$ws.Range("J80").Value = 25199.5
$ws.Range("L80").Value = 75598.5
$ws.Range("N80").Value = -77594.5
$ws.Range("H83").Value = 25199.5
$ws.Range("J83").Value = 25199.5
$ws.Range("L83").Value = 226795.5
$ws.Range("N83").Value = -236779.5
$ws.Range("H93").Value = 40000
$ws.Range("J93").Value = 40000
$ws.Range("L93").Value = 40000
$ws.Range("N93").Value = -44992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 5000
$ws.Range("J29").Value = 5000
$ws.Range("L29").Value = 5000
$ws.Range("N29").Value = -5616
$ws.Range("H32").Value = 2572.0833
$ws.Range("I32").Value = 2736.5
$ws.Range("J32").Value = 1750
$ws.Range("K32").Value = 2736.5
$ws.Range("L32").Value = 1750
$ws.Range("M32").Value = -2449.5
$ws.Range("N32").Value = -2324
$ws.Range("H45").Value = 1245
$ws.Range("I45").Value = 1245
$ws.Range("K45").Value = 1245
$ws.Range("M45").Value = -868
$ws.Range("H122").Value = 2253.3333
$ws.Range("I122").Value = 2030
$ws.Range("K122").Value = 6090
$ws.Range("M122").Value = -3640

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1833.3334
$ws.Range("H92").Value = 28000
$ws.Range("J92").Value = 28000
$ws.Range("L92").Value = 28000
$ws.Range("N92").Value = -32992
$ws.Range("H110").Value = 70348.5
$ws.Range("J110").Value = 70348.5
$ws.Range("L110").Value = 70348.5
$ws.Range("N110").Value = -78528.5
$ws.Range("H113").Value = 1833.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 75.52631
$ws.Range("I2").Value = 82.59999999999999
$ws.Range("J2").Value = 67.666664
$ws.Range("K2").Value = 82.59999999999999
$ws.Range("L2").Value = 67.666664
$ws.Range("M2").Value = 30.40000000000001
$ws.Range("N2").Value = -293.666664
$ws.Range("H15").Value = 25000
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H58").Value = 22500
$ws.Range("I58").Value = 25000
$ws.Range("J58").Value = 20000
$ws.Range("K58").Value = 25000
$ws.Range("L58").Value = 20000
$ws.Range("M58").Value = -24723
$ws.Range("N58").Value = -20554
$ws.Range("H70").Value = 1750
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 1500
$ws.Range("L70").Value = 2000
$ws.Range("M70").Value = -1230
$ws.Range("N70").Value = -2540
$ws.Range("H73").Value = 1750
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 1500
$ws.Range("L73").Value = 2000
$ws.Range("M73").Value = -564
$ws.Range("N73").Value = -3872
$ws.Range("H81").Value = 25000
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 25000
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0
$ws.Range("H113").Value = 1749.5
$ws.Range("J113").Value = 1999
$ws.Range("L113").Value = 1999
$ws.Range("N113").Value = -6339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 950
$ws.Range("I16").Value = 950
$ws.Range("K16").Value = 950
$ws.Range("M16").Value = -780
$ws.Range("H22").Value = 1744.8096
$ws.Range("I22").Value = 1506.7693
$ws.Range("J22").Value = 2131.625
$ws.Range("K22").Value = 1506.7693
$ws.Range("L22").Value = 2131.625
$ws.Range("M22").Value = -1211.7693
$ws.Range("N22").Value = -2721.625
$ws.Range("H27").Value = 1744.8096
$ws.Range("I27").Value = 1506.7693
$ws.Range("J27").Value = 2131.625
$ws.Range("K27").Value = 1506.7693
$ws.Range("L27").Value = 2131.625
$ws.Range("M27").Value = -1399.7693
$ws.Range("N27").Value = -2345.625
$ws.Range("H82").Value = 1945.5714
$ws.Range("J82").Value = 1809.7778
$ws.Range("L82").Value = 1809.7778
$ws.Range("N82").Value = -2531.7778
$ws.Range("H85").Value = 1945.5714
$ws.Range("J85").Value = 1809.7778
$ws.Range("L85").Value = 1809.7778
$ws.Range("N85").Value = -4305.7778
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("N94").Value = 0
$ws.Range("H132").Value = 7000
$ws.Range("I132").Value = 7000
$ws.Range("K132").Value = 21000
$ws.Range("M132").Value = -18470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2074.75
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 2099.6667
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 2099.6667
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -3347.6667
$ws.Range("H65").Value = 2074.75
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 2099.6667
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 10498.3335
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -16738.3335
$ws.Range("H100").Value = 756.25
$ws.Range("I100").Value = 785.7143
$ws.Range("J100").Value = 550
$ws.Range("K100").Value = 1571.4286
$ws.Range("L100").Value = 1100
$ws.Range("M100").Value = -1030.4286
$ws.Range("N100").Value = -2182
$ws.Range("H107").Value = 2549.5
$ws.Range("I107").Value = 799.5
$ws.Range("J107").Value = 4299.5
$ws.Range("K107").Value = 2398.5
$ws.Range("L107").Value = 12898.5
$ws.Range("M107").Value = -478.5
$ws.Range("N107").Value = -16738.5
$ws.Range("H122").Value = 2539.125
$ws.Range("I122").Value = 2539.125
$ws.Range("K122").Value = 7617.375
$ws.Range("M122").Value = -5167.375

